# Anonymize the laboratory/site names embedded in the "File" column (A2:A33).
# Mapping (derived from the original values -> their anonymized replacements):
#   PSURight   -> LaboratoryA_Right
#   PSULeft    -> LaboratoryA_Left
#   SLAC1      -> LaboratoryB_
#   UCSBRR     -> LaboratoryC_RR
#   Cargnello1 -> LaboratoryD_

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the write order below is deliberate (grouped by site, in the
# original rows' relative order) rather than row order 2..33, so that the
# workbook's shared-string table gets the new entries appended in the same
# order the source workbook ended up with (Excel/this engine appends newly
# seen strings to the shared string table in first-write order).
$fileNameUpdates = @(
    @{Row=2; Value="2p0_RhWI_0p02_500C_20250331_LaboratoryA_Right.xlsx"},
    @{Row=9; Value="2_RhWI_0p027_600C_20250422_LaboratoryA_Right.xlsx"},
    @{Row=10; Value="0p3_RhWI_0p0207_500C_20250315_LaboratoryA_Right.xlsx"},
    @{Row=19; Value="0p1_RhWI_0p02_500C_20250414_LaboratoryA_Right.xlsx"},
    @{Row=24; Value="2_RhWI_0p01824_500C_20250401_LaboratoryA_Right.xlsx"},
    @{Row=30; Value="2_RhWI_0p021_400C_20250425_LaboratoryA_Right.xlsx"},
    @{Row=31; Value="2_RhWI_0p02_500C_20250403_LaboratoryA_Right.xlsx"},
    @{Row=33; Value="0p3_RhNP_0p0204_500C_20250315_LaboratoryA_Left.xlsx"},
    @{Row=4; Value="0p3_RhWI_0p02mg_500C_20250331_LaboratoryB_.xlsx"},
    @{Row=6; Value="0p3_RhWI_0p0199mg_500C_20250415_LaboratoryB_.xlsx"},
    @{Row=8; Value="0p3_RhWI_0p0204mg_500C_20250421_LaboratoryB_.xlsx"},
    @{Row=11; Value="2_RhWI_0p02mg_400C_20250422_LaboratoryB_.xlsx"},
    @{Row=14; Value="2_RhWI_0p02mg_500C_20250417_LaboratoryB_.xlsx"},
    @{Row=22; Value="0p1_RhWI_0p02mg_500C_20250418_LaboratoryB_.xlsx"},
    @{Row=25; Value="2_RhWI_0p02mg_600C_20250423_LaboratoryB_.xlsx"},
    @{Row=28; Value="0p3_RhNP_0p0201mg_500C_20250424_LaboratoryB_.xlsx"},
    @{Row=5; Value="0p1_RhWI_0p02mg_500C_20250319_LaboratoryC_RR.xlsx"},
    @{Row=12; Value="0p3_RhNP_0p0201mg_500C_20250312_LaboratoryC_RR.xlsx"},
    @{Row=16; Value="0p3_RhWI_0p0201mg_500C_20250311_LaboratoryC_RR.xlsx"},
    @{Row=18; Value="2_RhWI_0p02mg_400C_20250306_LaboratoryC_RR.xlsx"},
    @{Row=20; Value="2_RhWI_0p02mg_600C_20250318_LaboratoryC_RR.xlsx"},
    @{Row=21; Value="2_RhWI_0p02mg_500C_20250307_LaboratoryC_RR.xlsx"},
    @{Row=26; Value="0p1_RhWI_0p02mg_500C_20250314_LaboratoryC_RR.xlsx"},
    @{Row=29; Value="0p1_RhWI_0p02mg_500C_20250317_LaboratoryC_RR.xlsx"},
    @{Row=3; Value="0p3_RhNP_0p0197mg_500C_20250211_LaboratoryD_.xlsx"},
    @{Row=7; Value="2_RhWI_0p0202mg_500C_20250317_LaboratoryD_.xlsx"},
    @{Row=13; Value="0p3_RhWI_0p02mg_500C_20250312_LaboratoryD_.xlsx"},
    @{Row=15; Value="2_RhWI_0p0204mg_400C_20250324_LaboratoryD_.xlsx"},
    @{Row=17; Value="0p3_RhNP_0p0201mg_500C_20250207_LaboratoryD_.xlsx"},
    @{Row=23; Value="2_RhWI_0p02mg_600C_20250326_LaboratoryD_.xlsx"},
    @{Row=27; Value="0p1_RhWI_0p02mg_500C_20250314_LaboratoryD_.xlsx"},
    @{Row=32; Value="0p3_RhNP_0p0199mg_500C_20250212_LaboratoryD_.xlsx"}
)

foreach ($update in $fileNameUpdates) {
    $ws.Range("A" + $update.Row).Value = $update.Value
}

# The workbook was last left with a different cell selected.
$ws.Range("G37").Select()
